# Fix bold-markdown typo in HLT_dict: a stray space before the closing
# "**" on "tones heard **" breaks the bold formatting. Remove that space
# ("tones heard**") in the three cells (C10, C11, C12) that contain it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = "<h4>How many tones can you hear?</h4>Your task is **to count again only the tones heard** and to enter the numerical value in the input box. **You should ignore the noise and not count it**.`n"
$ws.Range("C11").Value = "<h4>How many tones can you hear?</h4>Your task is **to count again only the tones heard** and to enter the numerical value in the input box. **You should ignore the noise and not count it**."
$ws.Range("C12").Value = "<h4>How many tones can you hear?</h4>Your task is **to count again only the tones heard** and to enter the numerical value in the input box. **You should ignore the noise and not count it**."
